$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 206, pushing existing rows 206..244 down to 207..245.
$ws.Rows("206:206").Insert()

# Populate the newly inserted row 206 with the new price-report entry.
$ws.Range("A206").Value = 4
$ws.Range("B206").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C206").Value = "Los Lagos"
$ws.Range("D206").Value = 44694
$ws.Range("E206").Value = 10
$ws.Range("F206").Value = 100112017
$ws.Range("G206").Value = "Apio"
$ws.Range("H206").Value = "Americana (o)"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 45
$ws.Range("K206").Value = 12000
$ws.Range("L206").Value = 12000
$ws.Range("M206").Value = 12000
$ws.Range("N206").Value = "`$/docena de matas"
$ws.Range("O206").Value = "Región de Coquimbo"
$ws.Range("P206").Value = 2000
$ws.Range("Q206").Value = 6
$ws.Range("R206").Value = "Hortaliza"
